# RPA datasets push 2024-07-03
# Insert two new rows (for 이노스페이스 IPO) right after the header row,
# pushing all existing data rows down by two, then populate the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank rows before the current row 2 (shifts rows 2-22 -> 4-24).
$ws.Range("A2:A3").EntireRow.Insert()

# 2) The inserted rows inherit the header row's formatting/border (style index).
#    Plain data rows in this sheet carry no explicit style, so strip it back off.
$ws.Range("A2:Q3").ClearFormats()

# 3) Force column A / O / P to Text so date-look-alike strings such as
#    "2024-07-02" are stored as literal text, not auto-converted to date serials.
$dateLikeCols = "A2","O2","P2","A3","O3","P3"
foreach ($addr in $dateLikeCols) {
    $ws.Range($addr).NumberFormat = "@"
}

# 4) Row 2: 이노스페이스 / 코스닥 / 미래 (대표)
$ws.Range("A2").Value = "2024-07-02"
$ws.Range("B2").Value = "이노스페이스"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 575.89
$ws.Range("E2").Value = "미래"
$ws.Range("F2").Value = 541.3366
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 43300
$ws.Range("N2").Value = 94
$ws.Range("O2").Value = "2024-06-20"
$ws.Range("P2").Value = "2024-06-25"
$ws.Range("Q2").Value = 1906400

# 5) Row 3: 이노스페이스 / 코스닥 / 신한 (인수)
$ws.Range("A3").Value = "2024-07-02"
$ws.Range("B3").Value = "이노스페이스"
$ws.Range("C3").Value = "코스닥"
$ws.Range("D3").Value = 575.89
$ws.Range("E3").Value = "신한"
$ws.Range("F3").Value = 34.5534
$ws.Range("G3").Value = "-"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = "-"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "인수"
$ws.Range("L3").Value = "-"
$ws.Range("M3").Value = 43300
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = "2024-06-20"
$ws.Range("P3").Value = "2024-06-25"
$ws.Range("Q3").Value = 1906400

# 6) Drop the text-format override again so the new rows end up with the
#    same "no explicit style" look as every other data row.
$ws.Range("A2:Q3").ClearFormats()
